# "Generate Report for Handback"
#
# The 99901ce2-f1a8-4724-9fbb-fe044b187417.md file has been handed back and
# is now in sync with en-US, so every "Ready for handoff" status for that
# file (in the Overview rollup as well as the per-locale zh-cn / de-de
# sheets) flips to "Handed back: in sync with en-US", and the per-locale
# handoff/handback timestamps + stale error detail are refreshed to reflect
# the new handback.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns for the handed-back file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: Status, Latest Handoff/Handback Datetime, Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-08-30 02:51:31"
$wsZhCn.Range("K3").Value = "2016-08-30 02:51:31"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: Status, Latest Handback Datetime, Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-30 02:51:38"
$wsDeDe.Range("P3").Value = ""
